# Plantilla Lista de Tareas de la Iteración - iteration update
# Sets "Horas estimadas totales" (G) and Día 9 consumption (AF) for the
# newly-worked tasks (rows 19 & 21), and records the Día 9 consumption
# for row 22 (which finishes off its remaining hours).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Fila 19 - "Realizar mockup de CU 14 y 16": se estima 1 hora y se
# consume completa en el día correspondiente (columna AF).
$ws.Range("G19").Value = 1
$ws.Range("AF19").Value = 1

# Fila 21 - "Realizar descripciones de CU 09 y 11": se estima 2 horas y
# se consumen completas (columna AF). Esta celda además recibe un
# formato que oculta el valor (texto y relleno en negro).
$ws.Range("G21").Value = 2
$ws.Range("AF21").Value = 2
$ws.Range("AF21").Font.ThemeColor = 1
$ws.Range("AF21").Interior.ThemeColor = 1

# Fila 22 - "Realizar descripciones de CU 10 y 12": se consume la hora
# restante (columna AF), agotando el saldo pendiente.
$ws.Range("AF22").Value = 1

# Re-merge the trailing "Total"/day-summary header cells so they end up
# at the tail of the merge list (matches how the workbook was re-saved).
$ws.Range("AZ4:BA4").UnMerge()
$ws.Range("AZ4:BA4").Merge()
$ws.Range("AO4:AP4").UnMerge()
$ws.Range("AO4:AP4").Merge()
$ws.Range("AR4:AS4").UnMerge()
$ws.Range("AR4:AS4").Merge()
$ws.Range("AU4:AV4").UnMerge()
$ws.Range("AU4:AV4").Merge()
$ws.Range("AX4:AY4").UnMerge()
$ws.Range("AX4:AY4").Merge()

$ws.Range("AC24").Select()
